# Reorder the output of Energy Star 5.2
#
# B12 (GPU Frame Buffer Width) is switched from ">= 128-bit" to "< 64-bit"
# (the dropdown list on B12 is "< 64-bit,>= 64-bit and < 128-bit,>= 128-bit"),
# and the two "Category B"/"Category C" formulas that used to fall back to the
# literal text "N/A" now fall back to an empty string instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pick a new value for the GPU Frame Buffer Width dropdown (B12).
$ws.Range("B12").Value = "< 64-bit"

# Column B (Category B) no longer reports "N/A" when the discrete-graphics
# test fails -- it reports a blank string instead.
$ws.Range("H2").Formula = '=IF(EXACT(B11,"Discrete"), "B", "")'

# Column C (Category C) gets the same treatment.
$ws.Range("I2").Formula = '=IF(AND(EXACT(B11,"Discrete"), EXACT(B12, ">= 128-bit"), B4>=2, B6>=2), "C", "")'

# Leave the selection on B11 (reflecting where the user clicked afterwards).
$ws.Range("B11").Select() | Out-Null
